$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 24 used to be an empty spacer row (only border/alignment formatting).
# It now becomes a fully populated "GET /dj/{dj-id}/songs" route row, matching
# the pattern of the other GET rows (e.g. row 9, row 15, row 20):
#   B = route, C = "GET", D = description, E = "application/json", F = response codes
# ---------------------------------------------------------------------------

# 1) Copy the cell formatting (borders/alignment) from the analogous GET row (row 9)
#    onto row 24, so the same style indices (border + vertical-center [+ wrap]) are reused.
$ws.Range("B9:F9").Copy() | Out-Null
$ws.Range("B24:F24").PasteSpecial(-4122) | Out-Null

# 2) Fill in the new text. "Alle Songs anzeigen" must be added to the shared
#    string table before "/dj/{dj-id}/songs" so they land on indices 37/38
#    respectively.
$ws.Range("D24").Value = "Alle Songs anzeigen"
$ws.Range("B24").Value = "/dj/{dj-id}/songs"
$ws.Range("C24").Value = "GET"
$ws.Range("E24").Value = "application/json"

# 3) F24 reuses the existing "200 OK / 404 not found / 500 internal error" text
#    (same as F9/F15/F20). Copy the source cell directly instead of retyping the
#    multi-line string so the existing shared string gets reused.
$ws.Range("F9").Copy() | Out-Null
$ws.Range("F24").PasteSpecial(-4104) | Out-Null

$excel.CutCopyMode = $false

# 4) Row 24 grows to fit the wrapped, multi-line content (same height as the
#    other "ht=45" description rows).
$ws.Rows.Item(24).RowHeight = 45

# ---------------------------------------------------------------------------
# The sheet's view was scrolled down and the selection moved from B33 to B24.
# ---------------------------------------------------------------------------
$ws.Range("B24").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
